# Updated cryptos list on Wed Jun 19 18:42:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.032.70'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.561.63'
$ws.Range("E3").Value = '  +4.36%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.78'
$ws.Range("E5").Value = '  +3.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.71'
$ws.Range("E6").Value = '  +3.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.560.86'
$ws.Range("E7").Value = '  +4.35%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +3.41%  '
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.169.58'
$ws.Range("E13").Value = '  +4.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000183'
$ws.Range("E14").Value = '  +3.36%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.574.35'
$ws.Range("E15").Value = '  +4.74%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.25'
$ws.Range("E16").Value = '  +4.85%  '
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.952.90'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.10'
$ws.Range("E19").Value = '  +7.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.40'
$ws.Range("E20").Value = '  +7.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.86'
$ws.Range("E21").Value = '  +3.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.40'
$ws.Range("E22").Value = '  +2.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.577'
$ws.Range("E23").Value = '  +7.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.709.14'
$ws.Range("E24").Value = '  +4.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.10'
$ws.Range("E25").Value = '  +3.41%  '
$ws.Range("E27").Value = '  +13.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.66'
$ws.Range("E28").Value = '  +6.45%  '
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.30'
$ws.Range("E30").Value = '  +5.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.32'
$ws.Range("E31").Value = '  +5.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.572.31'
$ws.Range("E32").Value = '  +4.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.43'
$ws.Range("E33").Value = '  +20.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.01'
$ws.Range("E34").Value = '  +5.25%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '170.60'
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.94'
$ws.Range("E38").Value = '  +5.18%  '
$ws.Range("E39").Value = '  +7.54%  '
$ws.Range("E40").Value = '  +9.76%  '
$ws.Range("E41").Value = '  +7.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.829'
$ws.Range("E42").Value = '  +3.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.96'
$ws.Range("E43").Value = '  +21.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.65'
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.46'
$ws.Range("E46").Value = '  +5.08%  '
$ws.Range("E47").Value = '  +10.11%  '
$ws.Range("E48").Value = '  +3.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.463.80'
$ws.Range("E49").Value = '  +12.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.90'
$ws.Range("E50").Value = '  +6.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.37'
$ws.Range("E51").Value = '  +17.06%  '
